# Refresh the cryptos list (prices / 1h volume %) with the latest scrape.
# Price column values that look like plain numbers are written with a
# leading apostrophe so Excel keeps them as text (matching the sheet's
# existing text-formatted "Price" column, e.g. "0.999", "396.80").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.937.13'
$ws.Range('E2').Value = '  +10.77%  '
$ws.Range('D3').Value = '3.260.63'
$ws.Range('E3').Value = '  +5.96%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''396.80'
$ws.Range('E5').Value = '  +0.73%  '
$ws.Range('D6').Value = '''109.50'
$ws.Range('E6').Value = '  +6.85%  '
$ws.Range('D7').Value = '''0.560'
$ws.Range('E7').Value = '  +4.81%  '
$ws.Range('D8').Value = '''0.999'
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = '''0.626'
$ws.Range('E9').Value = '  +6.50%  '
$ws.Range('D10').Value = '''39.31'
$ws.Range('E10').Value = '  +4.80%  '
$ws.Range('D11').Value = '''0.0958'
$ws.Range('E11').Value = '  +12.28%  '
$ws.Range('E12').Value = '  +2.31%  '
$ws.Range('D13').Value = '3.769.48'
$ws.Range('E13').Value = '  +5.99%  '
$ws.Range('D14').Value = '''8.24'
$ws.Range('E14').Value = '  +7.03%  '
$ws.Range('D15').Value = '''19.10'
$ws.Range('E15').Value = '  +2.59%  '
$ws.Range('D16').Value = '3.249.51'
$ws.Range('E16').Value = '  +6.87%  '
$ws.Range('E17').Value = '  +1.07%  '
$ws.Range('D18').Value = '''10.83'
$ws.Range('E18').Value = '  +2.62%  '
$ws.Range('D19').Value = '56.703.59'
$ws.Range('E19').Value = '  +10.36%  '
$ws.Range('D20').Value = '''3.30'
$ws.Range('E20').Value = '  +4.26%  '
$ws.Range('E21').Value = '  +9.65%  '
$ws.Range('D22').Value = '''12.89'
$ws.Range('E22').Value = '  +4.28%  '
$ws.Range('D23').Value = '''305.84'
$ws.Range('E23').Value = '  +15.40%  '
$ws.Range('D24').Value = '''75.11'
$ws.Range('E24').Value = '  +6.87%  '
$ws.Range('D25').Value = '''3.15'
$ws.Range('E25').Value = '  -1.82%  '
$ws.Range('D26').Value = '''28.13'
$ws.Range('E26').Value = '  +4.19%  '
$ws.Range('E27').Value = '  +4.89%  '
$ws.Range('D28').Value = '''7.89'
$ws.Range('E28').Value = '  +0.14%  '
$ws.Range('E29').Value = '  +2.60%  '
$ws.Range('E30').Value = '  +0.57%  '
$ws.Range('E31').Value = '  -0.38%  '
$ws.Range('D32').Value = '''0.110'
$ws.Range('E32').Value = '  +4.83%  '
$ws.Range('D33').Value = '''11.02'
$ws.Range('E33').Value = '  +2.71%  '
$ws.Range('D34').Value = '''37.48'
$ws.Range('E34').Value = '  +2.75%  '
$ws.Range('D35').Value = '''0.0480'
$ws.Range('E35').Value = '  -2.51%  '
$ws.Range('D36').Value = '''2.14'
$ws.Range('E36').Value = '  +3.19%  '
$ws.Range('D37').Value = '''51.50'
$ws.Range('E37').Value = '  +3.18%  '
$ws.Range('E38').Value = '  +5.51%  '

# Rows 39-40 swapped order (Stacks now ranks above FirstDigitalUSD).
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = '''3.09'
$ws.Range('E39').Value = '  +22.45%  '
$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').Value = '''0.998'
$ws.Range('E40').Value = '  -0.22%  '
$ws.Range('D41').Value = '''135.22'
$ws.Range('E41').Value = '  +4.87%  '
$ws.Range('E42').Value = '  +4.27%  '

# Rows 43-45 rotated (NEARProtocol, Stellar, Celestia shift up one rank).
$ws.Range('B43').Value = 'NEARProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D43').Value = '''4.01'
$ws.Range('E43').Value = '  -0.12%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').Value = '''0.120'
$ws.Range('E44').Value = '  +4.15%  '
$ws.Range('B45').Value = 'Celestia'
$ws.Range('C45').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D45').Value = '''17.17'
$ws.Range('E45').Value = '  +2.88%  '
$ws.Range('D46').Value = '''0.281'
$ws.Range('E46').Value = '  -2.59%  '
$ws.Range('D47').Value = '''22.06'
$ws.Range('E47').Value = '  +1.39%  '
$ws.Range('D48').Value = '2.147.73'
$ws.Range('E48').Value = '  +3.64%  '
$ws.Range('E49').Value = '  +1.98%  '
$ws.Range('E50').Value = '  -5.99%  '
$ws.Range('D51').Value = '''2.01'
$ws.Range('E51').Value = '  +37.35%  '
